$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 ("Data_resolution") previously had no Alias; give it one.
$ws.Range("D4").Value = "Data Resolution"

# Row 18 ("Year" field) renamed to "Survey_Year" / "Survey Year".
$ws.Range("A18").Value = "Survey_Year"
$ws.Range("D18").Value = "Survey Year"

# Move the active selection to E4, matching the saved cursor position.
$ws.Range("E4").Select()
